$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "SP23092022121400"
$ws.Range("B21").Value = 100
$ws.Range("C21").Value = "PANADOL STRIP 10"
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 15

$ws.Range("A22").Value = "SP23092022121400"
$ws.Range("B22").Value = 101
$ws.Range("C22").Value = "PANADOL STRIP 20"
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 29
